$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 38, shifting existing rows 38-40 down to 39-41
$ws.Rows.Item(38).Insert()

# Populate the new row 38 with the new weekly entry data
$ws.Cells.Item(38, 1).Value = 10
$ws.Cells.Item(38, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(38, 3).Value = "La Araucanía"
$ws.Cells.Item(38, 4).Value = 44585
$ws.Cells.Item(38, 4).NumberFormat = $ws.Cells.Item(39, 4).NumberFormat
$ws.Cells.Item(38, 5).Value = 9
$ws.Cells.Item(38, 6).Value = 100114002
$ws.Cells.Item(38, 7).Value = "Camote"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 20
$ws.Cells.Item(38, 11).Value = 16000
$ws.Cells.Item(38, 12).Value = 16000
$ws.Cells.Item(38, 13).Value = 16000
$ws.Cells.Item(38, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(38, 15).Value = "Perú"
$ws.Cells.Item(38, 16).Value = 800
$ws.Cells.Item(38, 17).Value = 20
$ws.Cells.Item(38, 18).Value = "Hortaliza"
